$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.75368041504542
$ws.Range("C2").Value = 10.24775696637774
$ws.Range("D2").Value = 5.964908691019467
$ws.Range("E2").Value = 16.55076751924465
$ws.Range("G2").Value = 3.601230190727962
$ws.Range("I2").Value = 17.95566946979838
$ws.Range("N2").Value = 15.91189494916883
$ws.Range("O2").Value = 19.98165726984459

$ws.Range("B3").Value = 13.08337459586296
$ws.Range("C3").Value = 9.596962367694028
$ws.Range("D3").Value = 5.842765516817798
$ws.Range("E3").Value = 15.60661971473264
$ws.Range("G3").Value = 3.603942032870336
$ws.Range("I3").Value = 18.03129993848473
$ws.Range("N3").Value = 15.96722713942471
$ws.Range("O3").Value = 19.93647452719252

$ws.Range("B4").Value = 12.65542097017515
$ws.Range("C4").Value = 9.172743985704713
$ws.Range("D4").Value = 5.768326594549031
$ws.Range("E4").Value = 15.00200625373128
$ws.Range("G4").Value = 3.605692935169996
$ws.Range("I4").Value = 18.08329142028107
$ws.Range("N4").Value = 16.0030372385761
$ws.Range("O4").Value = 19.91583893981606

$ws.Range("B5").Value = 12.47712839186676
$ws.Range("C5").Value = 8.993676589435038
$ws.Range("D5").Value = 5.738179336730623
$ws.Range("E5").Value = 14.74962803607787
$ws.Range("G5").Value = 3.606428101530311
$ws.Range("I5").Value = 18.10586615319201
$ws.Range("N5").Value = 16.01809295403244
$ws.Range("O5").Value = 19.90921766598145

$ws.Range("B6").Value = 12.44729479401609
$ws.Range("C6").Value = 8.963568440566895
$ws.Range("D6").Value = 5.73318609764854
$ws.Range("E6").Value = 14.70736774726427
$ws.Range("G6").Value = 3.606551485849305
$ws.Range("I6").Value = 18.1096982280942
$ws.Range("N6").Value = 16.02062093095619
$ws.Range("O6").Value = 19.90822617184665

$ws.Range("B7").Value = 12.65303191805154
$ws.Range("C7").Value = 9.170354103353695
$ws.Range("D7").Value = 5.767919197589298
$ws.Range("E7").Value = 14.99862646528523
$ws.Range("G7").Value = 3.605702762077623
$ws.Range("I7").Value = 18.08359026382248
$ws.Range("N7").Value = 16.00323840981146
$ws.Range("O7").Value = 19.91574240462558

$ws.Range("B8").Value = 13.52608164171616
$ws.Range("C8").Value = 10.0284720539351
$ws.Range("D8").Value = 5.922706050900596
$ws.Range("E8").Value = 16.23053943580562
$ws.Range("G8").Value = 3.602147469720422
$ws.Range("I8").Value = 17.98058883778998
$ws.Range("N8").Value = 15.93059313676776
$ws.Range("O8").Value = 19.9646044396161

$ws.Range("B9").Value = 15.1000664469454
$ws.Range("C9").Value = 11.51611301578408
$ws.Range("D9").Value = 6.228637760004749
$ws.Range("E9").Value = 18.544505278046
$ws.Range("G9").Value = 3.595852877254293
$ws.Range("I9").Value = 17.82306130162166
$ws.Range("N9").Value = 15.80265101021107
$ws.Range("O9").Value = 20.11666900011288

$ws.Range("B10").Value = 16.16351087716223
$ws.Range("C10").Value = 12.49072989441389
$ws.Range("D10").Value = 6.45223589913866
$ws.Range("E10").Value = 20.1992190691888
$ws.Range("G10").Value = 3.591636082005022
$ws.Range("I10").Value = 17.73495084049572
$ws.Range("N10").Value = 15.71742863312161
$ws.Range("O10").Value = 20.26233481341316

$ws.Range("B11").Value = 16.62574134741159
$ws.Range("C11").Value = 12.90857752158766
$ws.Range("D11").Value = 6.553179223102071
$ws.Range("E11").Value = 20.90982786788679
$ws.Range("G11").Value = 3.589805227706967
$ws.Range("I11").Value = 17.70097371003353
$ws.Range("N11").Value = 15.68054963580586
$ws.Range("O11").Value = 20.33585034083784

$ws.Range("B12").Value = 16.79759063485281
$ws.Range("C12").Value = 13.06315168793409
$ws.Range("D12").Value = 6.591251993350254
$ws.Range("E12").Value = 21.17289674979769
$ws.Range("G12").Value = 3.589124413943492
$ws.Range("I12").Value = 17.68899386296162
$ws.Range("N12").Value = 15.66685509352343
$ws.Range("O12").Value = 20.36471762165434

$ws.Range("B13").Value = 16.76072284402057
$ws.Range("C13").Value = 13.03002374911151
$ws.Range("D13").Value = 6.583059743941471
$ws.Range("E13").Value = 21.11650749245933
$ws.Range("G13").Value = 3.589270485019279
$ws.Range("I13").Value = 17.69153437183374
$ws.Range("N13").Value = 15.66979243152552
$ws.Range("O13").Value = 20.35845504045648

$ws.Range("B14").Value = 16.63994383186681
$ws.Range("C14").Value = 12.92136761984412
$ws.Range("D14").Value = 6.556314772078889
$ws.Range("E14").Value = 20.93159121924571
$ws.Range("G14").Value = 3.589748966839431
$ws.Range("I14").Value = 17.6999703040242
$ws.Range("N14").Value = 15.67941755680318
$ws.Range("O14").Value = 20.33820474485443

$ws.Range("B15").Value = 16.56554570367449
$ws.Range("C15").Value = 12.85433692262472
$ws.Range("D15").Value = 6.539911692095823
$ws.Range("E15").Value = 20.81754128787901
$ws.Range("G15").Value = 3.590043675245962
$ws.Range("I15").Value = 17.70525326544024
$ws.Range("N15").Value = 15.68534845480826
$ws.Range("O15").Value = 20.32593434262218

$ws.Range("B16").Value = 16.13286155315211
$ws.Range("C16").Value = 12.4629102965854
$ws.Range("D16").Value = 6.445619895348003
$ws.Range("E16").Value = 20.15193435525173
$ws.Range("G16").Value = 3.591757484659906
$ws.Range("I16").Value = 17.73729491852695
$ws.Range("N16").Value = 15.71987669672524
$ws.Range("O16").Value = 20.25767515406255

$ws.Range("B17").Value = 15.86183820984404
$ws.Range("C17").Value = 12.21625760606598
$ws.Range("D17").Value = 6.387547014333125
$ws.Range("E17").Value = 19.73284070766383
$ws.Range("G17").Value = 3.592831179577981
$ws.Range("I17").Value = 17.7585212780161
$ws.Range("N17").Value = 15.7415418352487
$ws.Range("O17").Value = 20.21764821054825

$ws.Range("B18").Value = 15.70392996206215
$ws.Range("C18").Value = 12.07198749960443
$ws.Range("D18").Value = 6.354074628348466
$ws.Range("E18").Value = 19.48782169402086
$ws.Range("G18").Value = 3.593456970248928
$ws.Range("I18").Value = 17.77130422272333
$ws.Range("N18").Value = 15.75418091283881
$ws.Range("O18").Value = 20.19530930158763

$ws.Range("B19").Value = 15.65012056138855
$ws.Range("C19").Value = 12.02272720402366
$ws.Range("D19").Value = 6.342730655243764
$ws.Range("E19").Value = 19.40417999171306
$ws.Range("G19").Value = 3.59367026797824
$ws.Range("I19").Value = 17.77573064188925
$ws.Range("N19").Value = 15.75849086381272
$ws.Range("O19").Value = 20.18786352442337

$ws.Range("B20").Value = 15.89089916689417
$ws.Range("C20").Value = 12.24276271508375
$ws.Range("D20").Value = 6.393736563619636
$ws.Range("E20").Value = 19.77786430226163
$ws.Range("G20").Value = 3.592716031721993
$ws.Range("I20").Value = 17.75620221462531
$ws.Range("N20").Value = 15.73921714293668
$ws.Range("O20").Value = 20.22183850038551

$ws.Range("B21").Value = 16.67550670029171
$ws.Range("C21").Value = 12.95338166623246
$ws.Range("D21").Value = 6.564174867740932
$ws.Range("E21").Value = 20.98606880384941
$ws.Range("G21").Value = 3.58960808672404
$ws.Range("I21").Value = 17.69746833735744
$ws.Range("N21").Value = 15.67658308370328
$ws.Range("O21").Value = 20.34412495700802

$ws.Range("B22").Value = 17.16968372006526
$ws.Range("C22").Value = 13.39651079079174
$ws.Range("D22").Value = 6.674661172123367
$ws.Range("E22").Value = 21.74060820534432
$ws.Range("G22").Value = 3.587649637619125
$ws.Range("I22").Value = 17.66425372231647
$ws.Range("N22").Value = 15.6372257444938
$ws.Range("O22").Value = 20.43003328464739

$ws.Range("B23").Value = 16.90766070888493
$ws.Range("C23").Value = 13.16194907342744
$ws.Range("D23").Value = 6.615788342813768
$ws.Range("E23").Value = 21.34109580164974
$ws.Range("G23").Value = 3.588688264804821
$ws.Range("I23").Value = 17.68150500091211
$ws.Range("N23").Value = 15.65808743772934
$ws.Range("O23").Value = 20.38363982341519

$ws.Range("B24").Value = 15.87776722138183
$ws.Range("C24").Value = 12.23078743796499
$ws.Range("D24").Value = 6.390938531605153
$ws.Range("E24").Value = 19.75752183556967
$ws.Range("G24").Value = 3.592768063556513
$ws.Range("I24").Value = 17.75724885832468
$ws.Range("N24").Value = 15.74026756465533
$ws.Range("O24").Value = 20.21994197257776

$ws.Range("B25").Value = 14.6900503550952
$ws.Range("C25").Value = 11.13450232586118
$ws.Range("D25").Value = 6.145896343363884
$ws.Range("E25").Value = 17.89755159066895
$ws.Range("G25").Value = 3.597483744731114
$ws.Range("I25").Value = 17.86086250821571
$ws.Range("N25").Value = 15.83571636818698
$ws.Range("O25").Value = 20.069536312289

Write-Output "Updated loading_percent values for 380 kV case"